$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update "actual time" (Фактические часы работы) for the week of 10-14 Dec ---
# Tuesday (E/F columns): leave time moved from 19:30 to 19:00, actual hours 6 -> 5.5
$ws.Range("E15").Value = 0.79166666666666663
$ws.Range("F14").Value = 5.5

# Wednesday (G/H columns): was not filled in yet, now filled with come/leave time and hours
$ws.Range("G14").Value = 0.67708333333333337
$ws.Range("H14").Value = 3.5
$ws.Range("G15").Value = 0.82291666666666663

# Thursday (I/J columns): was not filled in yet, now filled with come/leave time and hours
$ws.Range("I14").Value = 0.41666666666666669
$ws.Range("J14").Value = 3
$ws.Range("I15").Value = 0.54166666666666663

# New reviewer comment on H14 explaining the missed English class
$cmt = $ws.Range("H14").AddComment()
$cmt.Text("Anna Sharuntsova.EXT:" + [char]10 + "не пошла на английский")
$cmt.Shape.TextFrame.Characters(1, 21).Font.Bold = $true

$ws.Range("J30").Select()
